# Weekly fruit/vegetable price update:
# Insert a new daily-price record as row 178 on the "Ciboulette" sheet,
# pushing the existing rows 178-228 down to 179-229.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 178..228 down to 179..229, leaving a blank row 178 behind
# (mirrors Excel's own "Insert Copied/Blank Row" behaviour, including
# carrying the D-column date style down into the new row).
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with the new record.
$ws.Range("A178").Value = 4
$ws.Range("B178").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C178").Value = "Los Lagos"
$ws.Range("D178").Value = 44736
$ws.Range("E178").Value = 10
$ws.Range("F178").Value = 100112039
$ws.Range("G178").Value = "Ciboulette"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 240
$ws.Range("K178").Value = 2500
$ws.Range("L178").Value = 2500
$ws.Range("M178").Value = 2500
$ws.Range("N178").Value = "$/docena de atados"
$ws.Range("O178").Value = "Región Metropolitana"
$ws.Range("P178").Value = 833
$ws.Range("Q178").Value = 3
$ws.Range("R178").Value = "Hortaliza"
